# Update vm_pu.xlsx results for Case_4_141 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2,3,4,5,6,9,10,11,12,13,14)  # B,C,D,E,F,I,J,K,L,M,N

$data = @{}
$data[2] = @(1.02,1.052867787106645,1.051307780880689,1.056548272993526,1.059078839405825,1.045883159265253,1.057888129132677,1.054059644323757,1.059285681782752,1.061809323692522,1.05939045161751)
$data[3] = @(1.02,1.054238729137631,1.052346632113406,1.057869784064461,1.060618637801931,1.046335032841949,1.058907678995128,1.054910533833024,1.060419561517191,1.063161444746181,1.06041144935771)
$data[4] = @(1.02,1.05512460841034,1.053017711234352,1.058724027507972,1.061614131565259,1.046625549534432,1.059565723540369,1.055459388201459,1.061151851910797,1.064035001035596,1.061070428401684)
$data[5] = @(1.02,1.055496747503756,1.053299566481769,1.059082949719965,1.062032438143967,1.046747235716046,1.059841969157209,1.055689716172576,1.061459374634704,1.064401925266407,1.061347066318985)
$data[6] = @(1.02,1.055559214770333,1.053346875694152,1.059143202647344,1.062102662180752,1.046767641205399,1.059888328889722,1.05572836528136,1.061510989701333,1.064463514933828,1.061393491887636)
$data[7] = @(1.02,1.055129582062299,1.053021478440569,1.058728824230884,1.06161972177199,1.046627177265199,1.059569416300019,1.055462467463625,1.06115596234425,1.064039905140606,1.061074126405476)
$data[8] = @(1.02,1.053331357679377,1.05165910095111,1.056995065331307,1.059599403270089,1.046036261655293,1.058233039514003,1.054347566206068,1.059669175414082,1.062266563754633,1.059735851811139)
$data[9] = @(1.02,1.050153124404812,1.04924963183797,1.053933139676605,1.056032471024391,1.044980536220985,1.055865181601936,1.052369584882535,1.057038305979196,1.059131055855479,1.05736463126924)
$data[10] = @(1.02,1.048027541124763,1.047637194616571,1.051886967625165,1.053649479907797,1.044266871183366,1.054277617521316,1.051041723700307,1.055276747974103,1.057033181906251,1.055774812665529)
$data[11] = @(1.02,1.047105457492745,1.046937489766427,1.050999725170541,1.052616330983343,1.043955481558028,1.053587994007704,1.050464514246199,1.054512098648723,1.056122908666316,1.055084209807427)
$data[12] = @(1.02,1.046762693478308,1.046677357215547,1.050669972115305,1.052232370376435,1.043839459267498,1.053331501857442,1.05024977275385,1.054227785504481,1.055784502947953,1.054827353408891)
$data[13] = @(1.02,1.046836229449225,1.046733167054408,1.050740714035804,1.052314740550762,1.043864362698658,1.053386535548028,1.050295850994623,1.05428878477451,1.055857105350243,1.054882465253631)
$data[14] = @(1.02,1.047077129878696,1.046915991866374,1.050972471625316,1.052584596869797,1.043945898443774,1.053566799155495,1.05044677061154,1.054488603149074,1.05609494189618,1.055062984856097)
$data[15] = @(1.02,1.047225521779239,1.047028605525424,1.051115239544269,1.052750837099667,1.043996087749703,1.053677820934207,1.050539711964019,1.054611679473819,1.056241442234681,1.055174164298469)
$data[16] = @(1.02,1.048088700542333,1.047683599564339,1.051945824428081,1.053718018562603,1.044287486979019,1.054323338825084,1.051079983710697,1.055327455094748,1.057093553515552,1.055820598898793)
$data[17] = @(1.02,1.048629691918417,1.04809405300613,1.052466493166995,1.054324352730323,1.04446963824579,1.054727663157948,1.05141828013013,1.055775934256188,1.057627551973973,1.056225497418586)
$data[18] = @(1.02,1.048945080245639,1.048333318262514,1.052770071537661,1.054677892498401,1.044575655698956,1.054963286991314,1.051615387171634,1.056037343242229,1.057938843343358,1.056461455864822)
$data[19] = @(1.02,1.049052592110279,1.048414877017644,1.052873563928854,1.054798419516155,1.044611766264305,1.055043592872667,1.051682559145101,1.056126446359134,1.058044955191279,1.056541875789738)
$data[20] = @(1.02,1.048571665553249,1.048050030297025,1.052410642681322,1.054259311770371,1.044450118771993,1.05468430491967,1.051382006450997,1.055727835482094,1.057570277758355,1.056182077606637)
$data[21] = @(1.02,1.047006197969611,1.046862160928056,1.050904230131462,1.052505136581399,1.043921898111075,1.053513725312153,1.050402337970592,1.054429769621436,1.056024913017211,1.055009835641809)
$data[22] = @(1.02,1.046020412292895,1.046113961563152,1.049955976666059,1.051401037931785,1.043587710412182,1.052775793133155,1.049784411141785,1.053611953141277,1.055051601083099,1.054270855514475)
$data[23] = @(1.02,1.046543141758715,1.046510724551212,1.05045877117742,1.0519864559897,1.043765067123148,1.053167170756925,1.050112174005571,1.054045653338899,1.055567733815294,1.054662788939362)
$data[24] = @(1.02,1.048597885672452,1.048069922713004,1.052435879478717,1.054288701365192,1.04445893948723,1.05470389729501,1.051398397630911,1.055749569807602,1.057596158069446,1.056201697805397)
$data[25] = @(1.02,1.050975939629518,1.049873600808072,1.054725560169048,1.056955465596627,1.045255192887139,1.056478895571912,1.052882547840169,1.057719775971471,1.05994296053391,1.057979216783447)

foreach ($r in $data.Keys) {
    $rowData = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $rowData[$i]
    }
}

Write-Output "Updated vm_pu values for 380 kV case ($($data.Keys.Count) rows)"